$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D23").Value = "US Core Practitioner Profile"
$ws.Range("F26").Value = "US Core Practitioner Profile"
$ws.Range("E33").Value = "US Core Practitioner Profile"

# Row 34: shift existing E34 value ("US Core Specimen Profile") into F34,
# and set E34 to the new "US Core Practitioner Profile" value.
$ws.Range("F34").Value = $ws.Range("E34").Value2
$ws.Range("E34").Value = "US Core Practitioner Profile"

$ws.Range("D52").Value = "US Core Practitioner Profile"
